# Atualizado por script em 05-11-2023 08:45
# Appends 3 new match rows (90, 91, 92) at the end of the Ecuador Liga Pro 2023
# sheet, mirroring the formatting used by the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 89
$newRows = @(90, 91, 92)

# Duplicate formatting (styles / number formats) from the last existing row
# into the new rows before writing values, so the new rows inherit the same
# look (bold/bordered index column, datetime-formatted match-date column...).
foreach ($r in $newRows) {
    $src = $ws.Range("A" + $lastRow + ":V" + $lastRow)
    $dst = $ws.Range("A" + $r + ":V" + $r)
    $src.Copy($dst)
}

$data = @(
    @{
        Row = 90
        A = 89
        B = "ecuador"
        C = "liga-pro"
        D = "2023"
        E = 45234.79166666666
        F = "Gualaceo"
        G = 2
        H = "Delfin"
        I = 2
        J = 2.55
        K = "29/10/2023 19:13"
        L = 2.6
        M = "04/11/2023 18:53"
        N = 3.11
        O = "29/10/2023 19:13"
        P = 3.03
        Q = "04/11/2023 18:53"
        R = 2.94
        S = "29/10/2023 19:13"
        T = 3.04
        U = "04/11/2023 18:53"
        V = "https://www.betexplorer.com/football/ecuador/liga-pro/gualaceo-delfin/0n18Cpro/"
    },
    @{
        Row = 91
        A = 90
        B = "ecuador"
        C = "liga-pro"
        D = "2023"
        E = 45234.89583333334
        F = "Guayaquil City"
        G = 0
        H = "Emelec"
        I = 0
        J = 3.54
        K = "31/10/2023 01:12"
        L = 4.41
        M = "04/11/2023 21:20"
        N = 3.36
        O = "31/10/2023 01:12"
        P = 3.35
        Q = "04/11/2023 21:20"
        R = 2.11
        S = "31/10/2023 01:12"
        T = 1.92
        U = "04/11/2023 21:20"
        V = "https://www.betexplorer.com/football/ecuador/liga-pro/guayaquil-city-emelec/foIPK8z5/"
    },
    @{
        Row = 92
        A = 91
        B = "ecuador"
        C = "liga-pro"
        D = "2023"
        E = 45235
        F = "LDU Quito"
        G = 2
        H = "Ind. del Valle"
        I = 0
        J = 2.01
        K = "02/11/2023 01:12"
        L = 2.56
        M = "04/11/2023 23:58"
        N = 3.5
        O = "02/11/2023 01:12"
        P = 3.25
        Q = "04/11/2023 23:58"
        R = 3.66
        S = "02/11/2023 01:12"
        T = 2.89
        U = "04/11/2023 23:58"
        V = "https://www.betexplorer.com/football/ecuador/liga-pro/ldu-quito-independiente-del-valle/2BLXIn5H/"
    }
)

foreach ($row in $data) {
    $r = $row.Row

    $ws.Range("A" + $r).Value = $row.A
    $ws.Range("B" + $r).Value = $row.B
    $ws.Range("C" + $r).Value = $row.C

    # "temporada" column holds a text value that looks numeric ("2023").
    # Writing it through .Value directly would be auto-coerced to a real
    # number by Excel, so force text interpretation and then drop the
    # explicit number format again to keep the default (unstyled) cell
    # format used by the rest of the column.
    $dCell = $ws.Range("D" + $r)
    $dCell.NumberFormatLocal = "@"
    $dCell.Value = $row.D
    $dCell.ClearFormats()

    $ws.Range("E" + $r).Value = $row.E
    $ws.Range("F" + $r).Value = $row.F
    $ws.Range("G" + $r).Value = $row.G
    $ws.Range("H" + $r).Value = $row.H
    $ws.Range("I" + $r).Value = $row.I
    $ws.Range("J" + $r).Value = $row.J
    $ws.Range("K" + $r).Value = $row.K
    $ws.Range("L" + $r).Value = $row.L
    $ws.Range("M" + $r).Value = $row.M
    $ws.Range("N" + $r).Value = $row.N
    $ws.Range("O" + $r).Value = $row.O
    $ws.Range("P" + $r).Value = $row.P
    $ws.Range("Q" + $r).Value = $row.Q
    $ws.Range("R" + $r).Value = $row.R
    $ws.Range("S" + $r).Value = $row.S
    $ws.Range("T" + $r).Value = $row.T
    $ws.Range("U" + $r).Value = $row.U
    $ws.Range("V" + $r).Value = $row.V
}

Write-Host "Added rows 90-92 to sheet"
